$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Oct30 - Nov3"
$ws.Range("A1").Value = "Test"
